$wb = $excel.ActiveWorkbook

# --- 1. Rename existing sheets ---
$wsPos = $wb.Worksheets.Item(1)
$wsNeg = $wb.Worksheets.Item(2)
$wsPos.Name = "Get-Positive"
$wsNeg.Name = "Get-Negative"

# --- 2. Add new sheet at the end ---
$wsPost = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsPost.Name = "Post-Positive"

# --- 3. Column widths (approximate bestFit widths from the original authoring) ---
$wsPost.Columns.Item(2).ColumnWidth = 45.7265625
$wsPost.Columns.Item(3).ColumnWidth = 32.81640625
$wsPost.Columns.Item(4).ColumnWidth = 11.7265625
$wsPost.Columns.Item(5).ColumnWidth = 10.81640625
$wsPost.Columns.Item(6).ColumnWidth = 10.81640625
$wsPost.Columns.Item(7).ColumnWidth = 10.81640625
$wsPost.Columns.Item(8).ColumnWidth = 10.81640625
$wsPost.Columns.Item(9).ColumnWidth = 10.81640625
$wsPost.Columns.Item(10).ColumnWidth = 14.90625
$wsPost.Columns.Item(11).ColumnWidth = 19
$wsPost.Columns.Item(12).ColumnWidth = 29.08984375

# --- 4. Formatting: reuse the existing header / hyperlink-column styles already
#        present in the workbook so no redundant style entries are created ---
$wsPos.Range("A1").Copy()
$wsPost.Range("A1:L1").PasteSpecial(-4122)

$wsPos.Range("C2").Copy()
$wsPost.Range("C2:C5").PasteSpecial(-4122)
$wsPost.Range("D2:D5").PasteSpecial(-4122)
$wsPost.Range("J2:J5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 5. Header row ---
$wsPost.Cells.Item(1,1).Value  = '${TC_No}'
$wsPost.Cells.Item(1,2).Value  = 'Test Case'
$wsPost.Cells.Item(1,3).Value  = '${baseUrl}'
$wsPost.Cells.Item(1,4).Value  = '${relativeUrl}'
$wsPost.Cells.Item(1,5).Value  = '${firstName}'
$wsPost.Cells.Item(1,6).Value  = '${lastName}'
$wsPost.Cells.Item(1,7).Value  = '${userName}'
$wsPost.Cells.Item(1,8).Value  = '${password}'
$wsPost.Cells.Item(1,9).Value  = '${httpCode}'
$wsPost.Cells.Item(1,10).Value = '${email}'
$wsPost.Cells.Item(1,11).Value = '${successCode}'
$wsPost.Cells.Item(1,12).Value = '${respMessage}'

# --- 6. Data rows ---
# Row 2 - TC_01 new customer registration
$wsPost.Cells.Item(2,1).Value  = 'TC_01'
$wsPost.Cells.Item(2,2).Value  = 'Validate successful registration for new customer'
$wsPost.Cells.Item(2,3).Value  = 'http://restapi.demoqa.com/customer'
$wsPost.Cells.Item(2,4).Value  = '/register'
$wsPost.Cells.Item(2,5).Value  = 'John'
$wsPost.Cells.Item(2,6).Value  = 'Doe'
$wsPost.Cells.Item(2,7).Value  = 'new user'
$wsPost.Cells.Item(2,8).Value  = 'JDoe1558'
$wsPost.Cells.Item(2,9).Value  = 201
$wsPost.Cells.Item(2,10).Value = 'random email'
$wsPost.Cells.Item(2,11).Value = 'OPERATION_SUCCESS'
$wsPost.Cells.Item(2,12).Value = 'Operation completed successfully'

# Row 3 - TC_02 new customer registration
$wsPost.Cells.Item(3,1).Value  = 'TC_02'
$wsPost.Cells.Item(3,2).Value  = 'Validate successful registration for new customer'
$wsPost.Cells.Item(3,3).Value  = 'http://restapi.demoqa.com/customer'
$wsPost.Cells.Item(3,4).Value  = '/register'
$wsPost.Cells.Item(3,5).Value  = 'Jane'
$wsPost.Cells.Item(3,6).Value  = 'Doe'
$wsPost.Cells.Item(3,7).Value  = 'new user'
$wsPost.Cells.Item(3,8).Value  = 'JDoe1558'
$wsPost.Cells.Item(3,9).Value  = 201
$wsPost.Cells.Item(3,10).Value = 'random email'
$wsPost.Cells.Item(3,11).Value = 'OPERATION_SUCCESS'
$wsPost.Cells.Item(3,12).Value = 'Operation completed successfully'

# Row 4 - TC_03 new customer registration
$wsPost.Cells.Item(4,1).Value  = 'TC_03'
$wsPost.Cells.Item(4,2).Value  = 'Validate successful registration for new customer'
$wsPost.Cells.Item(4,3).Value  = 'http://restapi.demoqa.com/customer'
$wsPost.Cells.Item(4,4).Value  = '/register'
$wsPost.Cells.Item(4,5).Value  = 'Johnny'
$wsPost.Cells.Item(4,6).Value  = 'Depp'
$wsPost.Cells.Item(4,7).Value  = 'new user'
$wsPost.Cells.Item(4,8).Value  = 'JDoe1558'
$wsPost.Cells.Item(4,9).Value  = 201
$wsPost.Cells.Item(4,10).Value = 'random email'
$wsPost.Cells.Item(4,11).Value = 'OPERATION_SUCCESS'
$wsPost.Cells.Item(4,12).Value = 'Operation completed successfully'

# Row 5 - TC_04 existing customer registration
$wsPost.Cells.Item(5,1).Value  = 'TC_04'
$wsPost.Cells.Item(5,2).Value  = 'Validate registration scenario for existing customer'
$wsPost.Cells.Item(5,3).Value  = 'http://restapi.demoqa.com/customer'
$wsPost.Cells.Item(5,4).Value  = '/register'
$wsPost.Cells.Item(5,5).Value  = 'John'
$wsPost.Cells.Item(5,6).Value  = 'Doe'
$wsPost.Cells.Item(5,7).Value  = 'jdoe1234'
$wsPost.Cells.Item(5,8).Value  = 'JDoe1558'
$wsPost.Cells.Item(5,9).Value  = 200
$wsPost.Cells.Item(5,10).Value = 'jdoe@gmail.com'
$wsPost.Cells.Item(5,11).Value = 'User already exists'
$wsPost.Cells.Item(5,12).Value = 'FAULT_USER_ALREADY_EXISTS'

# --- 7. Hyperlinks (row 5 only, mirroring the base-url / email hyperlink pattern
#        already used on the Get-Positive / Get-Negative sheets). Hyperlinks.Add
#        resets the target cell's style, so re-apply the hyperlink-column format
#        afterwards to keep C5/J5 consistent with the rest of their columns. ---
$wsPost.Hyperlinks.Add($wsPost.Range("C5"), "http://restapi.demoqa.com/")
$wsPost.Hyperlinks.Add($wsPost.Range("J5"), "mailto:jdoe@gmail.com")

$wsPos.Range("C2").Copy()
$wsPost.Range("C5").PasteSpecial(-4122)
$wsPost.Range("J5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 8. Page setup ---
$wsPost.PageSetup.PaperSize = 9
$wsPost.PageSetup.Orientation = 1

# --- 9. Active-cell selection on the new sheet ---
$wsPost.Range("E15").Select()

# --- 10. Selection changes on the existing sheets (per original diff) ---
$wsPos.Range("A1:XFD2").Select()
$wsNeg.Range("A1:XFD1").Select()

# --- 11. Make sure the new sheet ends up the active tab/view ---
$wsPost.Activate()
$wsPost.Range("E15").Select()

Write-Output "done"
